$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing weekly case counts (rows 2-51, weeks 1-50)
$ws.Range("B2").Value = 306
$ws.Range("B3").Value = 469
$ws.Range("B4").Value = 479
$ws.Range("B5").Value = 673
$ws.Range("B6").Value = 481
$ws.Range("B7").Value = 535
$ws.Range("B8").Value = 538
$ws.Range("B9").Value = 588
$ws.Range("B10").Value = 533
$ws.Range("B11").Value = 601
$ws.Range("B12").Value = 524
$ws.Range("B13").Value = 433
$ws.Range("B14").Value = 490
$ws.Range("B15").Value = 513
$ws.Range("B16").Value = 458
$ws.Range("B17").Value = 291
$ws.Range("B18").Value = 492
$ws.Range("B19").Value = 423
$ws.Range("B20").Value = 441
$ws.Range("B21").Value = 368
$ws.Range("B22").Value = 416
$ws.Range("B23").Value = 407
$ws.Range("B24").Value = 395
$ws.Range("B25").Value = 315
$ws.Range("B26").Value = 385
$ws.Range("B27").Value = 292
$ws.Range("B28").Value = 269
$ws.Range("B29").Value = 368
$ws.Range("B30").Value = 377
$ws.Range("B31").Value = 470
$ws.Range("B32").Value = 346
$ws.Range("B33").Value = 290
$ws.Range("B34").Value = 323
$ws.Range("B35").Value = 434
$ws.Range("B36").Value = 466
$ws.Range("B37").Value = 510
$ws.Range("B38").Value = 580
$ws.Range("B39").Value = 526
$ws.Range("B40").Value = 557
$ws.Range("B41").Value = 517
$ws.Range("B42").Value = 392
$ws.Range("B43").Value = 347
$ws.Range("B44").Value = 478
$ws.Range("B45").Value = 447
$ws.Range("B46").Value = 451
$ws.Range("B47").Value = 484
$ws.Range("B48").Value = 400
$ws.Range("B49").Value = 468
$ws.Range("B50").Value = 458
$ws.Range("B51").Value = 446

# Append new weeks 51-53 (rows 52-54) for semana 52 de 2025 update
$ws.Range("A52").Value = 51
$ws.Range("B52").Value = 620
$ws.Range("A53").Value = 52
$ws.Range("B53").Value = 366
$ws.Range("A54").Value = 53
$ws.Range("B54").Value = 81
